# Rename wording in log/message cells from "value" to "amount" for the
# "Using previous value:" / "Using current value:" phrases, across all
# worksheets in the workbook (these phrases live in the AD column of the
# "Income" and "fond privat" sheets).

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rowCount = $used.Rows.Count
    $colCount = $used.Columns.Count
    $startRow = $used.Row
    $startCol = $used.Column

    for ($r = 0; $r -lt $rowCount; $r++) {
        for ($c = 0; $c -lt $colCount; $c++) {
            $cell = $ws.Cells.Item($startRow + $r, $startCol + $c)
            $val = $cell.Value()
            if ($val -ne $null -and $val -is [string]) {
                if ($val -like "*Using previous value*" -or $val -like "*Using current value*") {
                    $newVal = $val -replace "Using previous value", "Using previous amount"
                    $newVal = $newVal -replace "Using current value", "Using current amount"
                    $cell.Value = $newVal
                }
            }
        }
    }
}
